$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current values for the rows involved before mutating anything
$b4 = $ws.Range("B4").Value()
$c4 = $ws.Range("C4").Value()
$b5 = $ws.Range("B5").Value()
$c5 = $ws.Range("C5").Value()
$b6 = $ws.Range("B6").Value()
$c6 = $ws.Range("C6").Value()
$b7 = $ws.Range("B7").Value()
$c7 = $ws.Range("C7").Value()

# Reverse the order of rows 4-7 (both label and data columns)
$ws.Range("B4").Value = $b7
$ws.Range("C4").Value = $c7
$ws.Range("B5").Value = $b6
$ws.Range("C5").Value = $c6
$ws.Range("B6").Value = $b5
$ws.Range("C6").Value = $c5
$ws.Range("B7").Value = $b4
$ws.Range("C7").Value = $c4

# Remove row 10 entirely (the "edge betweenness centrality" row)
$ws.Rows(10).Delete()
